# feat: add phone_number to import graduations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Số điện thoại" (phone number) column (I) ---------------------

# Header cell I1: clone the format of the neighboring "Căn cước công dân"
# header (H1) so the new header picks up the same font/border/alignment,
# then overwrite its text.
$headerSrc = $ws.Cells.Item(1, 8)
$headerDst = $ws.Cells.Item(1, 9)
$headerSrc.Copy($headerDst)
$headerDst.Value = "Số điện thoại"

# Data cell I2: clone the format of the neighboring "Căn cước công dân"
# value (H2) - this keeps the quote-prefixed, wrap-text, bordered style -
# then set the phone number. A leading apostrophe forces the numeric-looking
# string to be stored as text (matching how the ID-like H2 value is stored).
$valueSrc = $ws.Cells.Item(2, 8)
$valueDst = $ws.Cells.Item(2, 9)
$valueSrc.Copy($valueDst)
$valueDst.Value = "'0987654321"

# --- Column width adjustments -------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 30.42578125
$ws.Columns.Item(4).ColumnWidth = 14.42578125
$ws.Columns.Item(9).ColumnWidth = 11.42578125

# --- Update the saved active selection ----------------------------------
$ws.Range("G6").Select()

$wb.Save()
